$d = $word.ActiveDocument

# --- 1) "Wn que va a arr" + "eglar la " bullet -----------------------------
# These two runs used to be split by the (now stale) "_GoBack" bookmark.
# Re-writing the whole phrase merges them back into a single run and drops
# the bookmark that sat between them.
$d.Content.Find.Execute("que va a arreglar la", $false, $false, $false, `
    $false, $false, $true, 1, $false, "que va a arreglar la", 2)

# --- 2) "Antes del evento:" paragraph --------------------------------------
# The "_GoBack" bookmark (Word's "last edit" marker) now belongs here,
# wrapping "Antes de", and the paragraph gained a trailing space typed as
# its own run.
$target = $d.Content
$target.Find.Execute("Antes del evento:", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$paraStart = $target.Start

$bmRange = $d.Range($paraStart, $paraStart + 8)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tail = $d.Range($paraStart + 17, $paraStart + 17)
$tail.InsertAfter(" ")
